# Generate Report for Handoff
#
# The localization run finished and the status moved from "In Translation"
# to "Ready for handoff"; the two "xliff generated" timestamps were
# refreshed to the moment the handoff package was produced. Because the
# new status text is longer than the old one, the Status column on each
# sheet is re-sized to fit it.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# Columns: E = zh-cn status, F = de-de status, G = Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-07 06:29:30"

# --- zh-cn detail sheet ----------------------------------------------
# Columns: C = Status, H = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-07 06:29:17"

# --- de-de detail sheet ----------------------------------------------
# Columns: C = Status, H = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-07 06:29:30"

# --- Widen the Status columns to fit the new, longer text ------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
